$d = $word.ActiveDocument

function Find-ParagraphIndex($needle) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        if ($d.Paragraphs($i).Range.Text.Contains($needle)) {
            return $i
        }
    }
    return -1
}

# ---------------------------------------------------------------------------
# Change 1: "Louis: Create enemy prefabs with the necessary scripts" -
# collapse the separate " " run and "Create enemy prefabs..." run into a
# single run reading " Create enemy prefabs with the necessary scripts".
# ---------------------------------------------------------------------------
$idx1 = Find-ParagraphIndex("Create enemy prefabs with the necessary scripts")
$p1 = $d.Paragraphs($idx1)
$full1 = $p1.Range
$text1 = $full1.Text
$needle1 = " Create enemy prefabs with the necessary scripts"
$pos1 = $full1.Start + $text1.IndexOf($needle1)

$spaceRange = $d.Range($pos1, $pos1 + 1)
$restRange = $d.Range($pos1 + 1, $full1.End - 1)
$restText = $restRange.Text
$restRange.Text = ""
$spaceRange2 = $d.Range($pos1, $pos1 + 1)
$spaceRange2.InsertAfter($restText)

# ---------------------------------------------------------------------------
# Change 2: "Khalid: Finish UI elements" -> "Khalid: Create UI elements",
# split across two runs: "Create" and " UI elements".
# ---------------------------------------------------------------------------
$idx2 = Find-ParagraphIndex("Finish UI elements")
$p2 = $d.Paragraphs($idx2)
$full2 = $p2.Range
$text2 = $full2.Text
$pos2 = $full2.Start + $text2.IndexOf("Finish")
$finishRange = $d.Range($pos2, $pos2 + 6)
# Force a run split by toggling character formatting, then restoring it -
# this keeps "Create" and " UI elements" as distinct runs afterwards.
$finishRange.Bold = 1
$finishRange.Text = "Create"
$createRange = $d.Range($pos2, $pos2 + 6)
$createRange.Bold = 0

# ---------------------------------------------------------------------------
# Change 3: insert three additional blank paragraphs - two right before the
# "At this point the game has modular mechanics..." paragraph, and one right
# after it.
# ---------------------------------------------------------------------------
$idx3 = Find-ParagraphIndex("At this point the game has modular mechanics")
$beforeP = $d.Paragraphs($idx3 - 1)
$beforeP.Range.InsertParagraphAfter()
$newP1 = $d.Paragraphs($idx3)
$newP1.Range.Delete()

$beforeP2 = $d.Paragraphs($idx3)
$beforeP2.Range.InsertParagraphAfter()
$newP2 = $d.Paragraphs($idx3 + 1)
$newP2.Range.Delete()

$idx3b = Find-ParagraphIndex("At this point the game has modular mechanics")
$targetP = $d.Paragraphs($idx3b)
$targetP.Range.InsertParagraphAfter()
$newP3 = $d.Paragraphs($idx3b + 1)
$newP3.Range.Delete()

Write-Output "done"
